$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44497
$ws.Range("J2").Value = 800
$ws.Range("K2").Value = 7500
$ws.Range("L2").Value = 8000
$ws.Range("M2").Value = 7750
$ws.Range("N2").Value = '$/caja 60 unidades'
$ws.Range("P2").Value = 129
$ws.Range("Q2").Value = 60
$ws.Range("D3").Value = 44321
$ws.Range("J3").Value = 500
$ws.Range("K3").Value = 7000
$ws.Range("L3").Value = 8000
$ws.Range("M3").Value = 7500
$ws.Range("N3").Value = '$/caja 60 unidades'
$ws.Range("P3").Value = 125
$ws.Range("Q3").Value = 60
$ws.Range("D4").Value = 44265
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 6500
$ws.Range("L4").Value = 7000
$ws.Range("M4").Value = 6750
$ws.Range("N4").Value = '$/caja 60 unidades'
$ws.Range("P4").Value = 112
$ws.Range("Q4").Value = 60
$ws.Range("D5").Value = 44258
$ws.Range("J5").Value = 500
$ws.Range("K5").Value = 7000
$ws.Range("L5").Value = 8000
$ws.Range("M5").Value = 7500
$ws.Range("N5").Value = '$/caja 60 unidades'
$ws.Range("P5").Value = 125
$ws.Range("Q5").Value = 60
$ws.Range("D6").Value = 44203
$ws.Range("J6").Value = 300
$ws.Range("K6").Value = 4500
$ws.Range("L6").Value = 5000
$ws.Range("M6").Value = 4750
$ws.Range("N6").Value = '$/caja 60 unidades'
$ws.Range("P6").Value = 79
$ws.Range("Q6").Value = 60
$ws.Range("D7").Value = 44300
$ws.Range("J7").Value = 400
$ws.Range("K7").Value = 6000
$ws.Range("L7").Value = 7000
$ws.Range("M7").Value = 6500
$ws.Range("N7").Value = '$/caja 60 unidades'
$ws.Range("P7").Value = 108
$ws.Range("Q7").Value = 60
$ws.Range("D8").Value = 44483
$ws.Range("J8").Value = 340
$ws.Range("K8").Value = 10000
$ws.Range("L8").Value = 11000
$ws.Range("M8").Value = 10500
$ws.Range("N8").Value = '$/caja 60 unidades'
$ws.Range("P8").Value = 175
$ws.Range("Q8").Value = 60
$ws.Range("D9").Value = 44251
$ws.Range("J9").Value = 700
$ws.Range("K9").Value = 6500
$ws.Range("L9").Value = 7000
$ws.Range("M9").Value = 6750
$ws.Range("N9").Value = '$/caja 60 unidades'
$ws.Range("P9").Value = 112
$ws.Range("Q9").Value = 60
$ws.Range("D10").Value = 44336
$ws.Range("J10").Value = 600
$ws.Range("K10").Value = 8500
$ws.Range("L10").Value = 9000
$ws.Range("M10").Value = 8750
$ws.Range("N10").Value = '$/caja 50 unidades'
$ws.Range("P10").Value = 175
$ws.Range("Q10").Value = 50
$ws.Range("D11").Value = 44266
$ws.Range("J11").Value = 600
$ws.Range("K11").Value = 6500
$ws.Range("L11").Value = 7000
$ws.Range("M11").Value = 6750
$ws.Range("N11").Value = '$/caja 60 unidades'
$ws.Range("P11").Value = 112
$ws.Range("Q11").Value = 60
$ws.Range("D12").Value = 44238
$ws.Range("J12").Value = 400
$ws.Range("K12").Value = 7000
$ws.Range("L12").Value = 8000
$ws.Range("M12").Value = 7500
$ws.Range("N12").Value = '$/caja 60 unidades'
$ws.Range("P12").Value = 125
$ws.Range("Q12").Value = 60
$ws.Range("D13").Value = 44377
$ws.Range("J13").Value = 400
$ws.Range("K13").Value = 7000
$ws.Range("L13").Value = 8000
$ws.Range("M13").Value = 7500
$ws.Range("N13").Value = '$/caja 60 unidades'
$ws.Range("P13").Value = 125
$ws.Range("Q13").Value = 60
$ws.Range("D14").Value = 44335
$ws.Range("J14").Value = 500
$ws.Range("K14").Value = 7500
$ws.Range("L14").Value = 8000
$ws.Range("M14").Value = 7750
$ws.Range("N14").Value = '$/caja 50 unidades'
$ws.Range("P14").Value = 155
$ws.Range("Q14").Value = 50
$ws.Range("D15").Value = 44216
$ws.Range("J15").Value = 1100
$ws.Range("K15").Value = 5500
$ws.Range("L15").Value = 6000
$ws.Range("M15").Value = 5750
$ws.Range("N15").Value = '$/caja 60 unidades'
$ws.Range("P15").Value = 96
$ws.Range("Q15").Value = 60
$ws.Range("D16").Value = 44314
$ws.Range("J16").Value = 1100
$ws.Range("K16").Value = 7000
$ws.Range("L16").Value = 8000
$ws.Range("M16").Value = 7500
$ws.Range("N16").Value = '$/caja 60 unidades'
$ws.Range("P16").Value = 125
$ws.Range("Q16").Value = 60
$ws.Range("D17").Value = 44294
$ws.Range("J17").Value = 500
$ws.Range("K17").Value = 7000
$ws.Range("L17").Value = 8000
$ws.Range("M17").Value = 7500
$ws.Range("N17").Value = '$/caja 60 unidades'
$ws.Range("P17").Value = 125
$ws.Range("Q17").Value = 60
$ws.Range("D18").Value = 44308
$ws.Range("J18").Value = 400
$ws.Range("K18").Value = 6000
$ws.Range("L18").Value = 7000
$ws.Range("M18").Value = 6500
$ws.Range("N18").Value = '$/caja 60 unidades'
$ws.Range("P18").Value = 108
$ws.Range("Q18").Value = 60
$ws.Range("D19").Value = 44315
$ws.Range("J19").Value = 500
$ws.Range("K19").Value = 7000
$ws.Range("L19").Value = 8000
$ws.Range("M19").Value = 7500
$ws.Range("N19").Value = '$/caja 60 unidades'
$ws.Range("P19").Value = 125
$ws.Range("Q19").Value = 60
$ws.Range("D20").Value = 44371
$ws.Range("J20").Value = 300
$ws.Range("K20").Value = 8500
$ws.Range("L20").Value = 9000
$ws.Range("M20").Value = 8750
$ws.Range("N20").Value = '$/caja 60 unidades'
$ws.Range("P20").Value = 146
$ws.Range("Q20").Value = 60
$ws.Range("D21").Value = 44490
$ws.Range("J21").Value = 600
$ws.Range("K21").Value = 13000
$ws.Range("L21").Value = 15000
$ws.Range("M21").Value = 14000
$ws.Range("N21").Value = '$/caja 60 unidades'
$ws.Range("P21").Value = 233
$ws.Range("Q21").Value = 60
$ws.Range("D22").Value = 44279
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 7000
$ws.Range("L22").Value = 8000
$ws.Range("M22").Value = 7500
$ws.Range("N22").Value = '$/caja 50 unidades'
$ws.Range("P22").Value = 150
$ws.Range("Q22").Value = 50
$ws.Range("D23").Value = 44482
$ws.Range("J23").Value = 400
$ws.Range("K23").Value = 11000
$ws.Range("L23").Value = 12000
$ws.Range("M23").Value = 11500
$ws.Range("N23").Value = '$/caja 60 unidades'
$ws.Range("P23").Value = 192
$ws.Range("Q23").Value = 60
$ws.Range("D24").Value = 44293
$ws.Range("J24").Value = 400
$ws.Range("K24").Value = 7000
$ws.Range("L24").Value = 8000
$ws.Range("M24").Value = 7500
$ws.Range("N24").Value = '$/caja 60 unidades'
$ws.Range("P24").Value = 125
$ws.Range("Q24").Value = 60
$ws.Range("D25").Value = 44384
$ws.Range("J25").Value = 300
$ws.Range("K25").Value = 7000
$ws.Range("L25").Value = 8000
$ws.Range("M25").Value = 7500
$ws.Range("N25").Value = '$/caja 60 unidades'
$ws.Range("P25").Value = 125
$ws.Range("Q25").Value = 60
$ws.Range("D26").Value = 44286
$ws.Range("J26").Value = 600
$ws.Range("K26").Value = 7000
$ws.Range("L26").Value = 8000
$ws.Range("M26").Value = 7500
$ws.Range("N26").Value = '$/caja 50 unidades'
$ws.Range("P26").Value = 150
$ws.Range("Q26").Value = 50
$ws.Range("D27").Value = 44181
$ws.Range("J27").Value = 900
$ws.Range("K27").Value = 4500
$ws.Range("L27").Value = 5000
$ws.Range("M27").Value = 4750
$ws.Range("N27").Value = '$/caja 60 unidades'
$ws.Range("P27").Value = 79
$ws.Range("Q27").Value = 60
$ws.Range("D28").Value = 44328
$ws.Range("J28").Value = 500
$ws.Range("K28").Value = 7500
$ws.Range("L28").Value = 8000
$ws.Range("M28").Value = 7750
$ws.Range("N28").Value = '$/caja 50 unidades'
$ws.Range("P28").Value = 155
$ws.Range("Q28").Value = 50
$ws.Range("D29").Value = 44301
$ws.Range("J29").Value = 300
$ws.Range("K29").Value = 6000
$ws.Range("L29").Value = 7000
$ws.Range("M29").Value = 6500
$ws.Range("N29").Value = '$/caja 60 unidades'
$ws.Range("P29").Value = 108
$ws.Range("Q29").Value = 60
$ws.Range("D30").Value = 44217
$ws.Range("J30").Value = 700
$ws.Range("K30").Value = 6500
$ws.Range("L30").Value = 7000
$ws.Range("M30").Value = 6750
$ws.Range("N30").Value = '$/caja 60 unidades'
$ws.Range("P30").Value = 112
$ws.Range("Q30").Value = 60
$ws.Range("D31").Value = 44244
$ws.Range("J31").Value = 500
$ws.Range("K31").Value = 5000
$ws.Range("L31").Value = 6000
$ws.Range("M31").Value = 5500
$ws.Range("N31").Value = '$/caja 60 unidades'
$ws.Range("P31").Value = 92
$ws.Range("Q31").Value = 60
$ws.Range("D32").Value = 44307
$ws.Range("J32").Value = 700
$ws.Range("K32").Value = 6000
$ws.Range("L32").Value = 7000
$ws.Range("M32").Value = 6500
$ws.Range("N32").Value = '$/caja 60 unidades'
$ws.Range("P32").Value = 108
$ws.Range("Q32").Value = 60
